# Updated cryptos list - applies price/volume/coin-order changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "46.093.41"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.598.74"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.16"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.39"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.29%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.584"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +2.57%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.12"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0844"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.54%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.09"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.17"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.993.98"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.31%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.594.30"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.918"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "46.217.99"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.09%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.33%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.77"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +2.18%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.84"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "289.55"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +14.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "73.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.07"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +3.40%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "29.25"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +5.21%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.90"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "39.26"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.48%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.04%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.27"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.45"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0840"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -4.41%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.88%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.97%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0333"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.11%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Celestia"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "15.71"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.59"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.77%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +10.85%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.109.22"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "97.36"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.999"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.48"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.86%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.203"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "109.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.05%  "
